$wb = $excel.ActiveWorkbook

# xlEdge* constants for Borders.Item(...)
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
# xlPasteFormats constant for PasteSpecial
$xlPasteFormats = -4122

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Derive the two new "box header" border-only styles exactly once, on
# quality_comparison!C1/D1, then fan them out to every other cell that needs
# the identical look via Copy/PasteSpecial(Formats). (Re-deriving the same
# look independently from scratch on each cell, or stripping edges back off
# a 4-sided border, tends to leave stray/orphan style or border records
# behind in this engine, so the look is built once - one border edge at a
# time, always landing on a combination the sheet's style table already
# knows about - and then cloned onto the rest of the cells.)

# C1: thin top+bottom border only (no left/right).
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item($xlEdgeTop).LineStyle = 1
$c1.Borders.Item($xlEdgeBottom).LineStyle = 1

# D1: thin top+right+bottom border (no left).
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item($xlEdgeTop).LineStyle = 1
$d1.Borders.Item($xlEdgeRight).LineStyle = 1
$d1.Borders.Item($xlEdgeBottom).LineStyle = 1

# Fan the two looks out to computational_comparison!C1/D1/F1/G1.
$c1.Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws2.Range("F1").PasteSpecial($xlPasteFormats)

$d1.Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell G5.
$ws2.Range("G5").ClearContents()
